function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# PIR sheet: add rows 144-157 (motion sensor, "No Motion" / "Inactive")
$wsPIR = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(144, "16:49:03"),
    @(145, "16:49:04"),
    @(146, "16:49:07"),
    @(147, "16:49:12"),
    @(148, "16:49:17"),
    @(149, "16:49:22"),
    @(150, "16:49:27"),
    @(151, "16:49:32"),
    @(152, "16:49:37"),
    @(153, "16:49:43"),
    @(154, "16:49:47"),
    @(155, "16:49:53"),
    @(156, "16:49:58"),
    @(157, "16:50:03")
)
foreach ($r in $pirRows) {
    $rowNum = $r[0]
    $ts = $r[1]
    Set-TextCell $wsPIR $rowNum 1 "2026-01-28"
    Set-TextCell $wsPIR $rowNum 2 $ts
    Set-TextCell $wsPIR $rowNum 3 "16:00"
    Set-TextCell $wsPIR $rowNum 4 "Bathroom"
    Set-TextCell $wsPIR $rowNum 5 "No Motion"
    Set-TextCell $wsPIR $rowNum 6 "Inactive"
}

# Humidity sheet: add rows 144-155 (humidity %, "Active")
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(144, "16:49:04", "87.9%"),
    @(145, "16:49:06", "87.9%"),
    @(146, "16:49:10", "87.0%"),
    @(147, "16:49:14", "87.9%"),
    @(148, "16:49:18", "86.9%"),
    @(149, "16:49:26", "87.8%"),
    @(150, "16:49:30", "86.9%"),
    @(151, "16:49:34", "87.9%"),
    @(152, "16:49:42", "86.9%"),
    @(153, "16:49:50", "86.9%"),
    @(154, "16:49:54", "87.8%"),
    @(155, "16:50:02", "87.9%")
)
foreach ($r in $humidityRows) {
    $rowNum = $r[0]
    $ts = $r[1]
    $val = $r[2]
    Set-TextCell $wsHumidity $rowNum 1 "2026-01-28"
    Set-TextCell $wsHumidity $rowNum 2 $ts
    Set-TextCell $wsHumidity $rowNum 3 "16:00"
    Set-TextCell $wsHumidity $rowNum 4 "Bathroom"
    Set-TextCell $wsHumidity $rowNum 5 $val
    Set-TextCell $wsHumidity $rowNum 6 "Active"
}

# Temperature sheet: add rows 144-155 (temperature C, "Active")
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @(144, "16:49:04", "22.8C"),
    @(145, "16:49:06", "22.8C"),
    @(146, "16:49:10", "22.8C"),
    @(147, "16:49:14", "22.8C"),
    @(148, "16:49:18", "22.8C"),
    @(149, "16:49:26", "22.8C"),
    @(150, "16:49:30", "22.8C"),
    @(151, "16:49:34", "22.8C"),
    @(152, "16:49:42", "22.8C"),
    @(153, "16:49:50", "22.8C"),
    @(154, "16:49:54", "22.8C"),
    @(155, "16:50:02", "22.9C")
)
foreach ($r in $temperatureRows) {
    $rowNum = $r[0]
    $ts = $r[1]
    $val = $r[2]
    Set-TextCell $wsTemperature $rowNum 1 "2026-01-28"
    Set-TextCell $wsTemperature $rowNum 2 $ts
    Set-TextCell $wsTemperature $rowNum 3 "16:00"
    Set-TextCell $wsTemperature $rowNum 4 "Bathroom"
    Set-TextCell $wsTemperature $rowNum 5 $val
    Set-TextCell $wsTemperature $rowNum 6 "Active"
}

Write-Host "Edit complete"
